# Update cryptocurrency price (D) and 1h volume/change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.407.23"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "2.318.21"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'518.87"
$ws.Range("E5").Value = "  +2.59%  "
$ws.Range("D6").Value = "'133.87"
$ws.Range("E6").Value = "  +3.69%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").Value = "2.337.39"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("E10").Value = "  +4.96%  "
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("D12").Value = "'5.28"
$ws.Range("E12").Value = "  +3.44%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "'23.78"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "2.732.79"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").Value = "56.531.99"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").Value = "2.357.43"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").Value = "'10.44"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "'322.91"
$ws.Range("E21").Value = "  +3.19%  "
$ws.Range("D22").Value = "'6.56"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'60.70"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("E25").Value = "  +5.44%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "'7.94"
$ws.Range("E27").Value = "  +6.01%  "
$ws.Range("E28").Value = "  +11.42%  "
$ws.Range("D29").Value = "0.0₃0737"
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("D31").Value = "'166.14"
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").Value = "'18.32"
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D35").Value = "'0.993"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").Value = "'0.920"
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("E39").Value = "  +6.38%  "
$ws.Range("D40").Value = "'37.81"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("D41").Value = "'0.382"
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("D42").Value = "'139.49"
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("D43").Value = "'3.59"
$ws.Range("E43").Value = "  +4.42%  "
$ws.Range("D44").Value = "'279.89"
$ws.Range("E44").Value = "  +7.44%  "
$ws.Range("D45").Value = "'5.17"
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("D46").Value = "'0.0931"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D49").Value = "'0.0217"
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").Value = "'17.75"
$ws.Range("E51").Value = "  +7.97%  "
